# Update the multiplication-fact answers in the single table of the document.
# Each data row of the table (rows 1,5,10,15,20) holds 5 cells like "A×B=C";
# we overwrite the w:t text of each target cell in place, preserving all
# paragraph/run formatting since only Range.Text is reassigned.
$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)
$x = [char]215   # the "×" multiplication sign used throughout the table

$cell = $tbl.Cell(1, 1)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "702${x}9=6318"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (1,1): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "887${x}9=7983"

$cell = $tbl.Cell(1, 2)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "386${x}9=3474"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (1,2): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "757${x}7=5299"

$cell = $tbl.Cell(1, 3)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "722${x}9=6498"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (1,3): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "912${x}2=1824"

$cell = $tbl.Cell(1, 4)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "650${x}8=5200"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (1,4): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "705${x}6=4230"

$cell = $tbl.Cell(1, 5)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "719${x}9=6471"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (1,5): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "705${x}3=2115"

$cell = $tbl.Cell(5, 1)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "487${x}9=4383"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (5,1): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "701${x}5=3505"

$cell = $tbl.Cell(5, 2)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "186${x}2=372"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (5,2): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "789${x}3=2367"

$cell = $tbl.Cell(5, 3)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "804${x}7=5628"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (5,3): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "792${x}3=2376"

$cell = $tbl.Cell(5, 4)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "609${x}6=3654"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (5,4): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "897${x}6=5382"

$cell = $tbl.Cell(5, 5)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "383${x}3=1149"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (5,5): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "880${x}5=4400"

$cell = $tbl.Cell(10, 1)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "636${x}2=1272"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (10,1): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "201${x}5=1005"

$cell = $tbl.Cell(10, 2)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "746${x}4=2984"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (10,2): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "850${x}8=6800"

$cell = $tbl.Cell(10, 3)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "710${x}8=5680"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (10,3): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "695${x}2=1390"

$cell = $tbl.Cell(10, 4)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "849${x}6=5094"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (10,4): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "462${x}5=2310"

$cell = $tbl.Cell(10, 5)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "367${x}5=1835"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (10,5): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "290${x}8=2320"

$cell = $tbl.Cell(15, 1)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "650${x}8=5200"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (15,1): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "936${x}4=3744"

$cell = $tbl.Cell(15, 2)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "652${x}9=5868"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (15,2): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "531${x}4=2124"

$cell = $tbl.Cell(15, 3)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "102${x}8=816"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (15,3): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "887${x}6=5322"

$cell = $tbl.Cell(15, 4)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "330${x}9=2970"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (15,4): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "713${x}2=1426"

$cell = $tbl.Cell(15, 5)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "654${x}2=1308"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (15,5): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "387${x}5=1935"

$cell = $tbl.Cell(20, 1)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "627${x}8=5016"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (20,1): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "959${x}9=8631"

$cell = $tbl.Cell(20, 2)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "623${x}9=5607"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (20,2): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "360${x}2=720"

$cell = $tbl.Cell(20, 3)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "401${x}5=2005"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (20,3): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "390${x}6=2340"

$cell = $tbl.Cell(20, 4)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "863${x}9=7767"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (20,4): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "456${x}5=2280"

$cell = $tbl.Cell(20, 5)
$actual = $cell.Range.Text.TrimEnd([char]13, [char]7)
$expectedOld = "221${x}3=663"
if ($actual -ne $expectedOld) {
  throw "Unexpected content in cell (20,5): expected [$expectedOld] but found [$actual]"
}
$cell.Range.Text = "963${x}6=5778"
